$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''29.035.47'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -0.86%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''1.830.82'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -0.87%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = '''0.9987'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''  +0.00%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''241.39'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +0.05%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''0.6539'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -2.94%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  +0.01%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''44.64'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  +5.89%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.07360'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -1.20%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''0.2939'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  -0.57%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = '''  -0.02%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''0.07675'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -0.55%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''1.829.81'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  -0.78%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''4.992'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -0.37%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''0.6672'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -0.87%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''82.36'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -4.55%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''6.071'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -1.34%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''0.000008630'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  +3.43%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''29.031.83'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -0.85%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''2.083.79'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -0.59%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''12.44'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -0.84%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''224.53'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -1.97%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''0.9997'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  -0.04%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''7.120'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  -1.32%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''1.000'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  +0.05%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''157.95'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -1.96%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''8.514'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -2.39%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''0.1382'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  -1.73%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''17.95'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  -0.42%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''1.501'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -0.68%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''4.110'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -1.69%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''1.206'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  +0.99%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''4.013'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  -1.44%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''0.05343'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  +0.58%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''0.7436'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  -2.16%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''1.833'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  -2.36%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''1.155'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  +1.39%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''2.642'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  -1.15%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''1.294.04'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  -2.35%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = '''  -1.14%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''2.746'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  +0.73%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''6.350'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  +6.07%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''0.8945'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  -2.60%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''0.9991'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = '''103.17'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -0.30%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''1.983.02'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -0.58%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''0.5141'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -0.52%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''64.17'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -0.23%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '''  -1.19%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''1.734'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -2.69%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''0.07522'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  -7.72%  '
$ws.Range("E51").Style = "Normal"
